$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1755102040816326
$ws.Range("C2").Value = 0.5265306122448979
$ws.Range("J2").Value = 0.00816326530612245
$ws.Range("O2").Value = 0.004081632653061225
$ws.Range("P2").Value = 0.1387755102040816
$ws.Range("S2").Value = 0.1469387755102041
$ws.Range("B3").Value = 0.0145985401459854
$ws.Range("C3").Value = 0.05109489051094891
$ws.Range("J3").Value = 0.04379562043795621
$ws.Range("P3").Value = 0.7518248175182481
$ws.Range("S3").Value = 0.1386861313868613
$ws.Range("J4").Value = 0.05
$ws.Range("P4").Value = 0.725
$ws.Range("S4").Value = 0.225
$ws.Range("B6").Value = 0.06862745098039216
$ws.Range("D6").Value = 0.0196078431372549
$ws.Range("F6").Value = 0.08333333333333333
$ws.Range("J6").Value = 0.196078431372549
$ws.Range("O6").Value = 0.0196078431372549
$ws.Range("Q6").Value = 0.2107843137254902
$ws.Range("R6").Value = 0.09313725490196079
$ws.Range("S6").Value = 0.3088235294117647
$ws.Range("B7").Value = 0.09547738693467336
$ws.Range("D7").Value = 0.03517587939698492
$ws.Range("E7").Value = 0.005025125628140704
$ws.Range("F7").Value = 0.06030150753768844
$ws.Range("J7").Value = 0.1105527638190955
$ws.Range("O7").Value = 0.01005025125628141
$ws.Range("Q7").Value = 0.1557788944723618
$ws.Range("R7").Value = 0.1206030150753769
$ws.Range("S7").Value = 0.407035175879397
$ws.Range("B8").Value = 0.06349206349206349
$ws.Range("D8").Value = 0.01587301587301587
$ws.Range("F8").Value = 0.04497354497354497
$ws.Range("J8").Value = 0.126984126984127
$ws.Range("O8").Value = 0.005291005291005291
$ws.Range("Q8").Value = 0.1693121693121693
$ws.Range("R8").Value = 0.1798941798941799
$ws.Range("S8").Value = 0.3941798941798942
$ws.Range("B9").Value = 0.1216931216931217
$ws.Range("D9").Value = 0.03703703703703703
$ws.Range("F9").Value = 0.04761904761904762
$ws.Range("J9").Value = 0.07407407407407407
$ws.Range("O9").Value = 0.01058201058201058
$ws.Range("Q9").Value = 0.1746031746031746
$ws.Range("R9").Value = 0.1111111111111111
$ws.Range("S9").Value = 0.4232804232804233
$ws.Range("B10").Value = 0.092
$ws.Range("D10").Value = 0.0144
$ws.Range("F10").Value = 0.0784
$ws.Range("J10").Value = 0.1152
$ws.Range("O10").Value = 0.0168
$ws.Range("Q10").Value = 0.2184
$ws.Range("R10").Value = 0.09520000000000001
$ws.Range("S10").Value = 0.3696
$ws.Range("G11").Value = 0.1375838926174497
$ws.Range("J11").Value = 0.08389261744966443
$ws.Range("K11").Value = 0.2147651006711409
$ws.Range("L11").Value = 0.5469798657718121
$ws.Range("S11").Value = 0.01677852348993289
$ws.Range("G12").Value = 0.7485029940119761
$ws.Range("J12").Value = 0.2155688622754491
$ws.Range("K12").Value = 0.01197604790419162
$ws.Range("L12").Value = 0.005988023952095809
$ws.Range("S12").Value = 0.01796407185628742
$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.3076923076923077
$ws.Range("F15").Value = 0.02094240837696335
$ws.Range("H15").Value = 0.1465968586387434
$ws.Range("I15").Value = 0.07329842931937172
$ws.Range("J15").Value = 0.3769633507853403
$ws.Range("K15").Value = 0.05759162303664921
$ws.Range("M15").Value = 0.01047120418848168
$ws.Range("O15").Value = 0.03664921465968586
$ws.Range("S15").Value = 0.2774869109947644
$ws.Range("H16").Value = 0.13125
$ws.Range("I16").Value = 0.11875
$ws.Range("J16").Value = 0.425
$ws.Range("K16").Value = 0.11875
$ws.Range("M16").Value = 0.03125
$ws.Range("O16").Value = 0.0625
$ws.Range("S16").Value = 0.1125
$ws.Range("F17").Value = 0.01809954751131222
$ws.Range("H17").Value = 0.167420814479638
$ws.Range("I17").Value = 0.09276018099547512
$ws.Range("J17").Value = 0.4638009049773756
$ws.Range("K17").Value = 0.07692307692307693
$ws.Range("M17").Value = 0.02036199095022624
$ws.Range("O17").Value = 0.05882352941176471
$ws.Range("S17").Value = 0.1018099547511312
$ws.Range("F18").Value = 0.01214574898785425
$ws.Range("H18").Value = 0.1336032388663968
$ws.Range("I18").Value = 0.08906882591093117
$ws.Range("J18").Value = 0.4736842105263158
$ws.Range("K18").Value = 0.1174089068825911
$ws.Range("M18").Value = 0.01619433198380567
$ws.Range("N18").Value = 0.004048582995951417
$ws.Range("O18").Value = 0.06072874493927125
$ws.Range("S18").Value = 0.0931174089068826
$ws.Range("F19").Value = 0.007718696397941681
$ws.Range("H19").Value = 0.1921097770154374
$ws.Range("I19").Value = 0.07890222984562607
$ws.Range("J19").Value = 0.3833619210977702
$ws.Range("K19").Value = 0.1157804459691252
$ws.Range("M19").Value = 0.02830188679245283
$ws.Range("N19").Value = 0.004288164665523156
$ws.Range("O19").Value = 0.06775300171526587
$ws.Range("S19").Value = 0.1217838765008576
